$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38
$ws.Cells.Item($row, 1).Value = "V-1770952993547"
$ws.Cells.Item($row, 2).Value = "12/2/2026"
$ws.Cells.Item($row, 3).Value = "10:23 p. m."
$ws.Cells.Item($row, 4).Value = "Stiven"
$ws.Cells.Item($row, 5).Value = "Aguardiente Amarillo Caja (x1)"
$ws.Cells.Item($row, 6).Value = 122000
$ws.Cells.Item($row, 7).Value = 0
